$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of extracted-sequence data appended below the existing ones
# (rows 13-16), matching the shape/format of the pre-existing rows.
$newRows = @(
    @{ Row = 13; A = "test"; B = 1741; C = 1851; D = 7; E = 1633; F = "33.28"; G = "35.38"; H = "0.13"; I = "31.21"; J = "2025-08-29 18:10:36" },
    @{ Row = 14; A = "test"; B = 1741; C = 1851; D = 7; E = 1633; F = "33.28"; G = "35.38"; H = "0.13"; I = "31.21"; J = "2025-08-29 18:13:53" },
    @{ Row = 15; A = "test"; B = 1741; C = 1851; D = 7; E = 1633; F = "33.28"; G = "35.38"; H = "0.13"; I = "31.21"; J = "2025-08-29 18:21:18" },
    @{ Row = 16; A = "test"; B = 1741; C = 1851; D = 7; E = 1633; F = "33.28"; G = "35.38"; H = "0.13"; I = "31.21"; J = "2025-08-29 18:23:34" }
)

# Columns F:J hold numeric/date-looking text (percentages and a timestamp)
# that must be stored as literal text, exactly like the existing rows.
# Pre-format that block as Text so the values aren't coerced to numbers.
$ws.Range("F13:J16").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
}
